$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - copy format from H1 (which already has the header style) to I1:J1
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I and J columns, rows 2-38
$data = @{
    2  = @(6,6)
    3  = @(7,7)
    4  = @(7,7)
    5  = @(8,8)
    6  = @(6,6)
    7  = @(8,8)
    8  = @(6,6)
    9  = @(8,9)
    10 = @(7,7)
    11 = @(7,7)
    12 = @(7,7)
    13 = @(6,6)
    14 = @(7,7)
    15 = @(5,6)
    16 = @(6,7)
    17 = @(6,7)
    18 = @(7,7)
    19 = @(9,9)
    20 = @(6,6)
    21 = @(7,7)
    22 = @(5,6)
    23 = @(9,9)
    24 = @(9,9)
    25 = @(9,9)
    26 = @(6,7)
    27 = @(6,7)
    28 = @(9,9)
    29 = @(8,8)
    30 = @(8,9)
    31 = @(8,8)
    32 = @(8,8)
    33 = @(7,7)
    34 = @(4,5)
    35 = @(7,7)
    36 = @(8,8)
    37 = @(5,6)
    38 = @(5,5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
